$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 6955.5454
$ws.Range("I43").Value = 9672.857
$ws.Range("J43").Value = 2200.25
$ws.Range("K43").Value = 9672.857
$ws.Range("L43").Value = 2200.25
$ws.Range("M43").Value = -9603.857
$ws.Range("N43").Value = -2338.25

$ws.Range("H51").Value = 4468.3335
$ws.Range("I51").Value = 3919.1667
$ws.Range("J51").Value = 5566.6665
$ws.Range("K51").Value = 3919.1667
$ws.Range("L51").Value = 5566.6665
$ws.Range("M51").Value = -3435.1667
$ws.Range("N51").Value = -6534.6665

$ws.Range("H98").Value = 6036.303
$ws.Range("I98").Value = 7046.154
$ws.Range("K98").Value = 7046.154
$ws.Range("M98").Value = -5548.154

$ws.Range("H116").Value = 15385.4
$ws.Range("I116").Value = 17227.5
$ws.Range("J116").Value = 14924.875
$ws.Range("K116").Value = 17227.5
$ws.Range("L116").Value = 14924.875
$ws.Range("M116").Value = -13785.5
$ws.Range("N116").Value = -21808.875

$ws.Range("H122").Value = 6036.303
$ws.Range("I122").Value = 7046.154
$ws.Range("K122").Value = 21138.462
$ws.Range("M122").Value = -18688.462

$ws.Range("H138").Value = 2492.8718
$ws.Range("I138").Value = 1871.25
$ws.Range("K138").Value = 5613.75
$ws.Range("M138").Value = -473.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4834.7607
$ws.Range("I32").Value = 2046.875
$ws.Range("K32").Value = 2046.875
$ws.Range("M32").Value = -1759.875

$ws.Range("H34").Value = 181404.4
$ws.Range("I34").Value = 101755.5
$ws.Range("J34").Value = 500000
$ws.Range("K34").Value = 101755.5
$ws.Range("L34").Value = 500000
$ws.Range("M34").Value = -101484.5
$ws.Range("N34").Value = -500542

$ws.Range("H110").Value = 21430708
$ws.Range("I110").Value = 34616068
$ws.Range("J110").Value = 4499.75
$ws.Range("K110").Value = 34616068
$ws.Range("L110").Value = 4499.75
$ws.Range("M110").Value = -34614023
$ws.Range("N110").Value = -8589.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2812.875

$ws.Range("H31").Value = 3896.3508
$ws.Range("I31").Value = 1978.826
$ws.Range("K31").Value = 1978.826
$ws.Range("M31").Value = -1683.826

$ws.Range("H34").Value = 3896.3508
$ws.Range("I34").Value = 1978.826
$ws.Range("K34").Value = 1978.826
$ws.Range("M34").Value = -1776.826

$ws.Range("H113").Value = 2812.875

$ws.Range("H122").Value = 2767.6667
$ws.Range("I122").Value = 2767
$ws.Range("J122").Value = 2769.8
$ws.Range("K122").Value = 8301
$ws.Range("L122").Value = 8309.400000000001
$ws.Range("M122").Value = -5851
$ws.Range("N122").Value = -13209.4

$ws.Range("H132").Value = 46029
$ws.Range("I132").Value = 3050.913
$ws.Range("K132").Value = 9152.739
$ws.Range("M132").Value = -6622.739

$ws.Range("H134").Value = 6415.0454
$ws.Range("I134").Value = 6029.05
$ws.Range("K134").Value = 18087.15
$ws.Range("M134").Value = -15552.15

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 72622.89999999999
$ws.Range("J37").Value = 72622.89999999999
$ws.Range("L37").Value = 217868.7
$ws.Range("N37").Value = -218092.7

$ws.Range("H108").Value = 3394.6667
$ws.Range("I108").Value = 3394.6667
$ws.Range("K108").Value = 10184.0001
$ws.Range("M108").Value = -7304.000100000001

$ws.Range("H109").Value = 3942.8572
$ws.Range("I109").Value = 600
$ws.Range("K109").Value = 1800
$ws.Range("M109").Value = -760

$ws.Range("H131").Value = 4973.222
$ws.Range("I131").Value = 709
$ws.Range("J131").Value = 7686.8184
$ws.Range("K131").Value = 2127
$ws.Range("L131").Value = 23060.4552
$ws.Range("M131").Value = 2913
$ws.Range("N131").Value = -33140.4552

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 500
$ws.Range("I29").Value = 500
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 500
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -210
$ws.Range("N29").ClearContents()

$ws.Range("H80").Value = 2164
$ws.Range("I80").Value = 1691.3334
$ws.Range("J80").Value = 5000
$ws.Range("K80").Value = 1691.3334
$ws.Range("L80").Value = 5000
$ws.Range("M80").Value = -693.3334
$ws.Range("N80").Value = -6996

$ws.Range("H83").Value = 2164
$ws.Range("I83").Value = 1691.3334
$ws.Range("J83").Value = 5000
$ws.Range("K83").Value = 8456.666999999999
$ws.Range("L83").Value = 25000
$ws.Range("M83").Value = -3464.666999999999
$ws.Range("N83").Value = -34984

$ws.Range("H126").Value = 3453.2727
$ws.Range("I126").Value = 2334
$ws.Range("K126").Value = 7002
$ws.Range("M126").Value = -4532

$ws.Range("H132").Value = 6211.5
$ws.Range("I132").Value = 3824.577
$ws.Range("J132").Value = 13107.056
$ws.Range("K132").Value = 11473.731
$ws.Range("L132").Value = 39321.16800000001
$ws.Range("M132").Value = -8943.731
$ws.Range("N132").Value = -44381.16800000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 4574.2856
$ws.Range("I61").Value = 5202.2
$ws.Range("K61").Value = 5202.2
$ws.Range("M61").Value = -5000.2

$ws.Range("H100").Value = 7356464
$ws.Range("I100").Value = 11907657
$ws.Range("J100").Value = 4536.4614
$ws.Range("K100").Value = 11907657
$ws.Range("L100").Value = 4536.4614
$ws.Range("M100").Value = -11907116
$ws.Range("N100").Value = -5618.4614

$ws.Range("H113").Value = 4574.2856
$ws.Range("I113").Value = 5202.2
$ws.Range("K113").Value = 5202.2
$ws.Range("M113").Value = -3032.2

$ws.Range("H122").Value = 7061.222
$ws.Range("I122").Value = 6679
$ws.Range("K122").Value = 20037
$ws.Range("M122").Value = -17587

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2084.8572
$ws.Range("J81").Value = 3773.2
$ws.Range("L81").Value = 7546.4
$ws.Range("N81").Value = -9668.4

$ws.Range("H84").Value = 2084.8572
$ws.Range("J84").Value = 3773.2
$ws.Range("L84").Value = 37732
$ws.Range("N84").Value = -48340

$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()

$ws.Range("H132").Value = 5254.865
$ws.Range("I132").Value = 3452.6177
$ws.Range("K132").Value = 10357.8531
$ws.Range("M132").Value = -7827.8531

$ws.Range("H136").Value = 5935.1
$ws.Range("I136").Value = 5349.3335
$ws.Range("K136").Value = 16048.0005
$ws.Range("M136").Value = -13498.0005
